$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 439, shifting existing rows 439+ down to 440+
$ws.Rows(439).Insert()

# Populate the new row 439 with data
$ws.Cells.Item(439, 1).Value = 10
$ws.Cells.Item(439, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(439, 3).Value = "La Araucanía"
$ws.Cells.Item(439, 4).Value = 44461
$ws.Cells.Item(439, 5).Value = 9
$ws.Cells.Item(439, 6).Value = "Fruta"
$ws.Cells.Item(439, 7).Value = 100102
$ws.Cells.Item(439, 8).Value = "Cítricos"
$ws.Cells.Item(439, 9).Value = 100102005
$ws.Cells.Item(439, 10).Value = "Naranja"
$ws.Cells.Item(439, 11).Value = "Navel Late"
$ws.Cells.Item(439, 12).Value = "Primera"
$ws.Cells.Item(439, 13).Value = 240
$ws.Cells.Item(439, 14).Value = 7000
$ws.Cells.Item(439, 15).Value = 8000
$ws.Cells.Item(439, 16).Value = 7417
$ws.Cells.Item(439, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(439, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(439, 19).Value = 494
$ws.Cells.Item(439, 20).Value = 15
